# Auto-generated Excel COM-interop edit script
# Applies the "Listado de reportes activos.xlsx" diff:
#  - Updates the COMMUNICATION (column E) text for specific rows: appends ",zip" to
#    some method lists and re-points a few cells at neighboring list variants
#    (shared-strings table is rebuilt automatically by Excel on save)
#  - Updates the active sheet selection (scroll is reset, selection moves to D9)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E620").Value = "correo, excel,factura"

$ws.Range("E897:E899").Value = "pdf,zip,correo,excel"
$ws.Range("E901:E903").Value = "pdf,zip,correo,excel"
$ws.Range("E948:E950").Value = "pdf,zip,correo,excel"
$ws.Range("E975:E981").Value = "pdf,zip,correo,excel"

$ws.Range("E220:E227").Value = "correo, excel, transfiere"

$ws.Range("E889").Value = "correo, ftp,pdf,excel"

$ws.Range("E995:E996").Value = "ftp,correo, pdf,excel"

$ws.Range("E1").Value = "COMUNICACIÓN`n(ej, Ftp, ws, smtp, etc,)"

$ws.Range("E997:E998").Value = "excel,correo,ftp,zip,xml, pdf"

$ws.Range("E189:E205").Value = "pdf,correo,excel,zip"
$ws.Range("E217:E219").Value = "pdf,correo,excel,zip"
$ws.Range("E228").Value = "pdf,correo,excel,zip"
$ws.Range("E255").Value = "pdf,correo,excel,zip"
$ws.Range("E257:E275").Value = "pdf,correo,excel,zip"
$ws.Range("E277").Value = "pdf,correo,excel,zip"
$ws.Range("E280:E314").Value = "pdf,correo,excel,zip"
$ws.Range("E316:E320").Value = "pdf,correo,excel,zip"
$ws.Range("E414:E443").Value = "pdf,correo,excel,zip"
$ws.Range("E458").Value = "pdf,correo,excel,zip"
$ws.Range("E508:E540").Value = "pdf,correo,excel,zip"
$ws.Range("E543:E561").Value = "pdf,correo,excel,zip"

$ws.Range("E594").Value = "ws,ftp,pdf,correo,excel,zip"
$ws.Range("E635:E641").Value = "ws,ftp,pdf,correo,excel,zip"
$ws.Range("E645:E658").Value = "ws,ftp,pdf,correo,excel,zip"
$ws.Range("E660:E665").Value = "ws,ftp,pdf,correo,excel,zip"
$ws.Range("E667:E671").Value = "ws,ftp,pdf,correo,excel,zip"
$ws.Range("E674:E693").Value = "ws,ftp,pdf,correo,excel,zip"
$ws.Range("E697:E700").Value = "ws,ftp,pdf,correo,excel,zip"
$ws.Range("E703").Value = "ws,ftp,pdf,correo,excel,zip"
$ws.Range("E708:E859").Value = "ws,ftp,pdf,correo,excel,zip"
$ws.Range("E878").Value = "ws,ftp,pdf,correo,excel,zip"
$ws.Range("E883").Value = "ws,ftp,pdf,correo,excel,zip"

# Update the sheet selection: this also clears the scrolled-away topLeftCell (was C1)
$ws.Range("D9").Select()
